# Update crypto price/volume snapshot cells to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as literal text in the workbook (values like
# "212.60" / "15.60" / "27.635.46" are not valid numbers to Excel because of the
# thousands-dot grouping and significant trailing zeros). Force the NumberFormat to
# Text ("@") before assigning so COM keeps them as strings instead of silently
# re-parsing them into doubles (which would drop the trailing zeros / grouping dots).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.635.46'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.636.80'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.60'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.524'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.97'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.633.38'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.558'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.51'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.616.07'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.90'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.74'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.05'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.48'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.96'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.60'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0486'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.453.11'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.563'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.897'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.91'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.778.19'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.35'

# Column E ("Volume(1h)") cells are percentage strings padded with spaces; they
# stay text naturally so a plain .Value assignment is enough.
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("E9").Value = '  -0.38%  '
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("E15").Value = '  -5.54%  '
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("E23").Value = '  +3.22%  '
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("E25").Value = '  +1.84%  '
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  +2.01%  '
$ws.Range("E34").Value = '  -2.44%  '
$ws.Range("E35").Value = '  -1.58%  '
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("E40").Value = '  +8.30%  '
$ws.Range("E41").Value = '  +7.81%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("E47").Value = '  -0.54%  '
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("E49").Value = '  -2.35%  '
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("E51").Value = '  -1.18%  '
